$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new "dcterms_subject" / Stichworte value for row 2 (column L)
$ws.Range("L2").Value = "Architektur; Personen"

# Update selection / view to reflect where the user ended up after editing
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("L3").Select()
